# "Generate Report for Handback"
# Updates the handback status report:
#  - flips the "in sync" status to "not in sync" everywhere it is shown
#  - records new Correspond-Handback timestamps for the f94e23ab file in
#    both the zh-cn and de-de sheets
#  - widens the "Status" columns to fit the longer status text

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: not in sync with en-US"

# --- Overview sheet: zh-cn (E) and de-de (F) status columns ---
$ws_overview.Range("E2").Value = $newStatus
$ws_overview.Range("F2").Value = $newStatus
$ws_overview.Range("E3").Value = $newStatus
$ws_overview.Range("F3").Value = $newStatus

# --- zh-cn detail sheet: Status column (C) ---
$ws_zhcn.Range("C2").Value = $newStatus
$ws_zhcn.Range("C3").Value = $newStatus

# --- de-de detail sheet: Status column (C) ---
$ws_dede.Range("C2").Value = $newStatus
$ws_dede.Range("C3").Value = $newStatus

# --- Updated handback timestamps for the f94e23ab file (row 3) ---
$ws_zhcn.Range("K3").Value = "2016-10-27 10:12:23"
$ws_dede.Range("K3").Value = "2016-10-27 10:12:39"

# --- Widen the Status columns so the longer text fits (report regenerated) ---
$newColumnWidth = 32.67
$ws_overview.Columns.Item(5).ColumnWidth = $newColumnWidth
$ws_overview.Columns.Item(6).ColumnWidth = $newColumnWidth
$ws_zhcn.Columns.Item(3).ColumnWidth = $newColumnWidth
$ws_dede.Columns.Item(3).ColumnWidth = $newColumnWidth
